$d = $word.ActiveDocument

# --- Change 1: "La nota de mi autoevaluación es 5.0, ..." -> "... es 5, ..." ---
$d.Content.Find.Execute(
    "La nota de mi autoevaluación es 5.0, ya que",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La nota de mi autoevaluación es 5, ya que", 2)

# --- Change 2: "La nota que le doy a mi compañero de trabajo es 5.0, ya que"
#     -> "La nota que le doy a mi compañero David es 5, ya que"
#     (ends up split across three runs, since "David" was inserted in the
#     middle of the original sentence and " de trabajo" was dropped while
#     "5.0" became "5") ---

# First drop " de trabajo" and fix "5.0" -> "5" (this keeps everything in one
# run for now, it will be split into two below).
$d.Content.Find.Execute(
    "compañero de trabajo es 5.0, ya que fue de mucho apoyo",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "compañero es 5, ya que fue de mucho apoyo", 2)

# Now insert " David" right after "compañero" and force it into its own run
# (matching the target's run boundaries) by nudging a character property on
# just the new text and reverting it - this splits the run without changing
# its visible formatting.
$rng = $d.Content
$rng.Find.Execute("La nota que le doy a mi compañero", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $rng.End
$insertPoint = $d.Range($s, $s)
$insertPoint.InsertAfter(" David")
$newRun = $d.Range($s, $s + 6)
$newRun.Bold = 1
$newRun.Bold = 0
